$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H7").Value = 5
$ws.Range("I7").Value = 5

$ws.Range("H4").Select()
